$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-92 down to 43-93.
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with data.
$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value = "Bíobío"
$ws.Cells.Item(42, 4).Value = 44539
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(42, 6).Value = 100112043
$ws.Cells.Item(42, 7).Value = "Pepino ensalada"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 180
$ws.Cells.Item(42, 11).Value = 6500
$ws.Cells.Item(42, 12).Value = 7000
$ws.Cells.Item(42, 13).Value = 6722
$ws.Cells.Item(42, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(42, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(42, 16).Value = 112
$ws.Cells.Item(42, 17).Value = 60
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D.
$ws.Cells.Item(42, 4).NumberFormat = $ws.Cells.Item(41, 4).NumberFormat
